$wb = $excel.ActiveWorkbook

# OFF sheet - Home row (row 2): update Short Att, Short Comp, Deep Att, Deep Comp
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 381
$wsOff.Range("C2").Value = 270
$wsOff.Range("D2").Value = 106
$wsOff.Range("E2").Value = 50

# DEF sheet - Home row (row 2): update Short Att, Short Comp, Deep Att, Deep Comp
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 454
$wsDef.Range("C2").Value = 300
$wsDef.Range("D2").Value = 91
$wsDef.Range("E2").Value = 35
